$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.022.79'
$ws.Range("E2").Value = '  +0.84%  '

$ws.Range("D3").Value = '2.471.35'
$ws.Range("E3").Value = '  +0.89%  '

$ws.Range("D5").Value = "'560.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.48%  '

$ws.Range("D6").Value = "'162.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = "'0.506"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.32%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = "'0.151"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.40%  '

$ws.Range("B10").Value = 'TRON'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D10").Value = "'0.165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.65%  '

$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").Value = "'0.332"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.56%  '

$ws.Range("D12").Value = "'4.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.60%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.905.82'
$ws.Range("E13").Value = '  +0.13%  '

$ws.Range("B14").Value = 'WrappedBTC'
$ws.Range("C14").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D14").Value = '68.899.57'
$ws.Range("E14").Value = '  +0.73%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = "'0.0000169"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.25%  '

$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Value = "'23.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.08%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.469.04'
$ws.Range("E17").Value = '  -2.31%  '

$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = "'10.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.22%  '

$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = "'336.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.64%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = "'6.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.82%  '

$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").Value = "'3.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.31%  '

$ws.Range("B22").Value = 'SuiNetwork'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D22").Value = "'1.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.30%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = "'66.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.77%  '

$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Value = '2.601.02'
$ws.Range("E25").Value = '  +1.15%  '

$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").Value = "'3.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.30%  '

$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.20%  '

$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").Value = "'8.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.51%  '

$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0816'
$ws.Range("E29").Value = '  -2.52%  '

$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = "'7.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.10%  '

$ws.Range("B31").Value = 'FirstDigitalUSD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.02%  '

$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").Value = "'430.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.38%  '

$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = "'1.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.26%  '

$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").Value = "'1.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.93%  '

$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").Value = "'158.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.03%  '

$ws.Range("B36").Value = 'WhiteBITCoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D36").Value = "'19.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.23%  '

$ws.Range("B37").Value = 'USDe'
$ws.Range("C37").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.05%  '

$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = "'0.109"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.03%  '

$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").Value = "'17.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.52%  '

$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").Value = "'0.300"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.01%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D41").Value = "'4.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.55%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = "'1.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.91%  '

$ws.Range("B43").Value = 'ImmutableX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D43").Value = "'1.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.34%  '

$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = "'2.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.94%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = "'132.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.99%  '

$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").Value = "'3.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.31%  '

$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = "'0.0713"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.72%  '

$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = "'0.483"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.12%  '

$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = "'0.562"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.12%  '

$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = "'0.0911"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.18%  '

$ws.Range("B51").Value = 'BitgetToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range("D51").Value = "'1.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.21%  '
